# Natmi following Dr Hou advice
# Update Lgi3-Adam23 LR-pair sheet: recompute rows 2-4 with new counts,
# and append two new target-cluster rows (M1, M2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs -> Lgi3/Adam23 -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Lgi3"
$ws.Range("C2").Value = "Adam23"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.607064
$ws.Range("H2").Value = 4.821192
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1790343333333333
$ws.Range("N2").Value = 0.537103
$ws.Range("O2").Value = 0.00485560647444994
$ws.Range("P2").Value = 0.00522711558998419
$ws.Range("Q2").Value = 0.287719631864
$ws.Range("R2").Value = 2.589476686776
$ws.Range("S2").Value = 0.00485560647444994
$ws.Range("T2").Value = 0.00522711558998419

# Row 3: FAPs -> Lgi3/Adam23 -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Lgi3"
$ws.Range("C3").Value = "Adam23"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.607064
$ws.Range("H3").Value = 4.821192
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 28.61718033333333
$ws.Range("N3").Value = 85.851541
$ws.Range("O3").Value = 0.7761291564580808
$ws.Range("P3").Value = 0.8355118634326505
$ws.Range("Q3").Value = 45.98964029520801
$ws.Range("R3").Value = 413.906762656872
$ws.Range("S3").Value = 0.7761291564580808
$ws.Range("T3").Value = 0.8355118634326505

# Row 4: FAPs -> Lgi3/Adam23 -> M1
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Lgi3"
$ws.Range("C4").Value = "Adam23"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.607064
$ws.Range("H4").Value = 4.821192
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.184753
$ws.Range("N4").Value = 0.5542590000000001
$ws.Range("O4").Value = 0.005010702954409397
$ws.Range("P4").Value = 0.005394078714490606
$ws.Range("Q4").Value = 0.2969098951920001
$ws.Range("R4").Value = 2.672189056728
$ws.Range("S4").Value = 0.005010702954409397
$ws.Range("T4").Value = 0.005394078714490606

# Row 5 (new): FAPs -> Lgi3/Adam23 -> M2
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Lgi3"
$ws.Range("C5").Value = "Adam23"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.607064
$ws.Range("H5").Value = 4.821192
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.02891466666666667
$ws.Range("N5").Value = 0.086744
$ws.Range("O5").Value = 0.0007841973104221829
$ws.Range("P5").Value = 0.0008441973229298452
$ws.Range("Q5").Value = 0.046467719872
$ws.Range("R5").Value = 0.418209478848
$ws.Range("S5").Value = 0.0007841973104221829
$ws.Range("T5").Value = 0.0008441973229298452

# Row 6 (new): FAPs -> Lgi3/Adam23 -> sCs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Lgi3"
$ws.Range("C6").Value = "Adam23"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.607064
$ws.Range("H6").Value = 4.821192
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 7.8617905
$ws.Range("N6").Value = 15.723581
$ws.Range("O6").Value = 0.2132203368026377
$ws.Range("P6").Value = 0.1530227449399449
$ws.Range("Q6").Value = 12.634400488092
$ws.Range("R6").Value = 75.80640292855199
$ws.Range("S6").Value = 0.2132203368026377
$ws.Range("T6").Value = 0.1530227449399449
